$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at V (index 22). Cell values/styles shift right
#    automatically. Comments stay anchored to their pre-insert cell refs in
#    this engine (they do not reflow with the column insert), so the rest of
#    this script reflows them explicitly to match what real Excel would do.
$ws.Columns("V:V").Insert()

# 2) New column V header text
$ws.Range("V15").Value2 = "culture_collection"

# 3) Reflow comments V15..CE15 -> W15..CF15. Processed highest column first so
#    each destination cell is read (via its own still-original comment object)
#    before it gets overwritten by the value shifting into it from the left.
#    This also preserves each comment's original author identity wherever a
#    comment object already exists at the destination.
$shiftData = @(
  @{ dest = "CF"; text = "specification of weight loss in the last three months, if yes should be further specified to include amount of weight loss" },
  @{ dest = "CE"; text = "history of urogenitaltract disorders; can include multiple disorders" },
  @{ dest = "CD"; text = "specification of urine collection method" },
  @{ dest = "CC"; text = "specification of twin sibling presence" },
  @{ dest = "CB"; text = "Feeding position in food chain (eg., chemolithotroph)" },
  @{ dest = "CA"; text = "specification of the countries travelled in the last six months; can include multiple travels" },
  @{ dest = "BZ"; text = "temperature of the sample at time of sampling" },
  @{ dest = "BY"; text = "Information about the genetic distinctness of the lineage (eg., biovar, serovar)" },
  @{ dest = "BX"; text = "specification of study completion status, if no the reason should be specified" },
  @{ dest = "BW"; text = "unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples." },
  @{ dest = "BV"; text = "specification of smoking status" },
  @{ dest = "BU"; text = "volume (mL) or weight (g) of sample processed for DNA extraction" },
  @{ dest = "BT"; text = "temperature at which sample was stored, e.g. -80" },
  @{ dest = "BS"; text = "location at which sample was stored, usually name of a specific freezer/room" },
  @{ dest = "BR"; text = "duration for which sample was stored" },
  @{ dest = "BQ"; text = "Amount or size of sample (volume, mass or area) that was collected" },
  @{ dest = "BP"; text = "salinity of sample, i.e. measure of total salt concentration" },
  @{ dest = "BO"; text = "Processing applied to the sample during or after isolation" },
  @{ dest = "BN"; text = "Method or device employed for collecting sample" },
  @{ dest = "BM"; text = "Aerobic or anaerobic" },
  @{ dest = "BL"; text = "history of pulmonary disorders; can include multiple disorders" },
  @{ dest = "BK"; text = "specification of presence of pets or farm animals in the environment of subject, if yes the animals should be specified; can include multiple animals present" },
  @{ dest = "BJ"; text = "type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types" },
  @{ dest = "BI"; text = "To what is the entity pathogenic" },
  @{ dest = "BH"; text = "oxygenation status of sample" },
  @{ dest = "BG"; text = "total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts" },
  @{ dest = "BF"; text = "history of nose-throat disorders; can include multiple disorders" },
  @{ dest = "BE"; text = "any other measurement performed or parameter collected, that is not listed here" },
  @{ dest = "BD"; text = "whether full medical history was collected" },
  @{ dest = "BC"; text = "specification of the maternal health status" },
  @{ dest = "BB"; text = "A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html" },
  @{ dest = "BA"; text = "history of kidney disorders; can include multiple disorders" },
  @{ dest = "AZ"; text = "Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived." },
  @{ dest = "AY"; text = "can include multiple medication codes" },
  @{ dest = "AX"; text = "total mass of the host at collection, the unit depends on host" },
  @{ dest = "AW"; text = "Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005" },
  @{ dest = "AV"; text = "NCBI taxonomy ID of the host, e.g. 9606" },
  @{ dest = "AU"; text = "a unique identifier by which each subject can be referred to, de-identified, e.g. #131" },
  @{ dest = "AT"; text = "Gender or physical sex of the host" },
  @{ dest = "AS"; text = "resting pulse of the host, measured as beats per minute" },
  @{ dest = "AQ"; text = "most frequent job performed by subject" },
  @{ dest = "AP"; text = "content of last meal and time since feeding; can include multiple values" },
  @{ dest = "AO"; text = "HIV status of subject, if yes HAART initiation status should also be indicated as [YES or NO]" },
  @{ dest = "AN"; text = "the height of subject" },
  @{ dest = "AK"; text = "Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh" },
  @{ dest = "AJ"; text = "type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types" },
  @{ dest = "AI"; text = "core body temperature of the host when sample was collected" },
  @{ dest = "AH"; text = "substance produced by the host, e.g. stool, mucus, where the sample was obtained from" },
  @{ dest = "AG"; text = "body mass index of the host, calculated as weight/(height)squared" },
  @{ dest = "AF"; text = "Age of host at the time of sampling" },
  @{ dest = "AE"; text = "Health or disease status of sample at time of collection" },
  @{ dest = "AD"; text = "specification of the gestation state" },
  @{ dest = "AC"; text = "specification of foetal health status, should also include abortion" },
  @{ dest = "AB"; text = "Plasmids that have significance phenotypic consequence" },
  @{ dest = "AA"; text = "ethnicity of the subject" },
  @{ dest = "Z"; text = "Estimated size of genome" },
  @{ dest = "Y"; text = "Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes" },
  @{ dest = "X"; text = "any drug used by subject and the frequency of usage; can include multiple drugs used" },
  @{ dest = "W"; text = "specification of major diet changes in the last six months, if yes the change should be specified" }
)
foreach ($item in $shiftData) {
    $destCell = $ws.Range($item.dest + "15")
    if ($destCell.Comment -ne $null) {
        $destCell.Comment.Text($item.text)
    } else {
        $destCell.AddComment($item.text)
    }
}

# 4) V15 keeps its own pre-existing comment object; give it the new
#    culture_collection text.
$vCell = $ws.Range("V15")
if ($vCell.Comment -ne $null) {
    $vCell.Comment.Text("Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier")
} else {
    $vCell.AddComment("Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier")
}

# 5) A couple of original comment-bearing columns (AM, AR) are themselves the
#    right-hand neighbour of a column that never had a comment (AL, AQ), so
#    nothing shifts into them in step 3 above and their stale original text
#    would otherwise be left behind. Remove those leftovers.
$orphanCols = @("AM", "AR")
foreach ($c in $orphanCols) {
    $cell = $ws.Range($c + "15")
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}
